$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# New error-code / text rows appended to the Translation sheet (rows 58-73).
# Columns: B=Text ID, C=Typography Name, D=Alignment, E=GB (text), F=Direction

$ws.Range("B58").Value = "SingleUseId76"
$ws.Range("C58").Value = "RpmText"
$ws.Range("D58").Value = "Center"
$ws.Range("E58").Value = "<value>"
$ws.Range("F58").Value = "LTR"
$ws.Range("B58:F58").Style = "Normal"

$ws.Range("B59").Value = "SingleUseId77"
$ws.Range("C59").Value = "RpmText"
$ws.Range("D59").Value = "Left"
$ws.Range("E59").NumberFormat = "@"
$ws.Range("E59").Value = "21"
$ws.Range("F59").Value = "LTR"
$ws.Range("B59:F59").Style = "Normal"

$ws.Range("B60").Value = "SingleUseId78"
$ws.Range("C60").Value = "Default"
$ws.Range("D60").Value = "Left"
$ws.Range("E60").Value = "IGN ANG"
$ws.Range("F60").Value = "LTR"
$ws.Range("B60:F60").Style = "Normal"

$ws.Range("B61").Value = "SingleUseId79"
$ws.Range("C61").Value = "Default"
$ws.Range("D61").Value = "Left"
$ws.Range("E61").Value = "<value>"
$ws.Range("F61").Value = "LTR"
$ws.Range("B61:F61").Style = "Normal"

$ws.Range("B62").Value = "SingleUseId80"
$ws.Range("C62").Value = "Default"
$ws.Range("D62").Value = "Left"
$ws.Range("E62").Value = "NO ERRORS"
$ws.Range("F62").Value = "LTR"
$ws.Range("B62:F62").Style = "Normal"

$ws.Range("B63").Value = "COOLANT_SENSOR_FAILURE"
$ws.Range("C63").Value = "Default"
$ws.Range("D63").Value = "Left"
$ws.Range("E63").Value = "Coolant temperature sensor failed!"
$ws.Range("F63").Value = "LTR"
$ws.Range("B63:F63").Style = "Normal"

$ws.Range("B64").Value = "NO_ERRORS"
$ws.Range("C64").Value = "Default"
$ws.Range("D64").Value = "Left"
$ws.Range("E64").Value = "NO ERRORS"
$ws.Range("F64").Value = "LTR"
$ws.Range("B64:F64").Style = "Normal"

$ws.Range("B65").Value = "KNOCKING"
$ws.Range("C65").Value = "Default"
$ws.Range("D65").Value = "Left"
$ws.Range("E65").Value = "Knock detected!"
$ws.Range("F65").Value = "LTR"
$ws.Range("B65:F65").Style = "Normal"

$ws.Range("B66").Value = "IAT_SENSOR_FAILURE"
$ws.Range("C66").Value = "Default"
$ws.Range("D66").Value = "Left"
$ws.Range("E66").Value = "IAT sensor failed!"
$ws.Range("F66").Value = "LTR"
$ws.Range("B66:F66").Style = "Normal"

$ws.Range("B67").Value = "MAP_SENSOR_FAILURE"
$ws.Range("C67").Value = "Default"
$ws.Range("D67").Value = "Left"
$ws.Range("E67").Value = "`nMAP sensor failed!"
$ws.Range("F67").Value = "LTR"
$ws.Range("B67:F67").Style = "Normal"

$ws.Range("B68").Value = "O2_SENSOR_FAILED"
$ws.Range("C68").Value = "Default"
$ws.Range("D68").Value = "Left"
$ws.Range("E68").Value = "O2 sensor failed!"
$ws.Range("F68").Value = "LTR"
$ws.Range("B68:F68").Style = "Normal"

$ws.Range("B69").Value = "EGT1_SESNSOR_FAILED"
$ws.Range("C69").Value = "Default"
$ws.Range("D69").Value = "Left"
$ws.Range("E69").Value = "EGT sensor #1 failed!"
$ws.Range("F69").Value = "LTR"
$ws.Range("B69:F69").Style = "Normal"

$ws.Range("B70").Value = "EGT2_SENSOR_FAILED"
$ws.Range("C70").Value = "Default"
$ws.Range("D70").Value = "Left"
$ws.Range("E70").Value = "EGT sensor #2 failed!"
$ws.Range("F70").Value = "LTR"
$ws.Range("B70:F70").Style = "Normal"

$ws.Range("B71").Value = "EGT_HIGH"
$ws.Range("C71").Value = "Default"
$ws.Range("D71").Value = "Left"
$ws.Range("E71").Value = "EGT too high!"
$ws.Range("F71").Value = "LTR"
$ws.Range("B71:F71").Style = "Normal"

$ws.Range("B72").Value = "DBW_SENSOR_FAILED"
$ws.Range("C72").Value = "Default"
$ws.Range("D72").Value = "Left"
$ws.Range("E72").Value = "Drive by wire failure!"
$ws.Range("F72").Value = "LTR"
$ws.Range("B72:F72").Style = "Normal"

$ws.Range("B73").Value = "FPR_RELATIVE_ERROR"
$ws.Range("C73").Value = "Default"
$ws.Range("D73").Value = "Left"
$ws.Range("E73").Value = "Fuel pressure relative error!"
$ws.Range("F73").Value = "LTR"
$ws.Range("B73:F73").Style = "Normal"
